# Generate Report for Handoff
# Updates "Latest Handoff Datetime" (column D) for rows whose files were just
# handed off, on both the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$rowsToUpdate = 4, 6, 7, 8, 9, 10

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $wsZhCn.Range("D$r").Value = "2016-03-08 02:54:02"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $wsDeDe.Range("D$r").Value = "2016-03-08 02:54:12"
}
